# simon.docx: append two trailing spaces to the first paragraph's existing
# text, then append a new, separately-formatted run in red (C00000) with the
# "(This is a change ... )" note, per the commit "changed simon and tiger".

$d = $word.ActiveDocument

# --- Step 1: extend "This is a Microsoft word document." with two spaces,
#     keeping it in the same (unformatted) run. Insert right before the
#     paragraph mark so it stays inside paragraph 1 rather than paragraph 2.
$para1 = $d.Paragraphs(1).Range
$textEnd = $para1.End - 1
$insertPoint = $d.Range($textEnd, $textEnd)
$insertPoint.InsertAfter("  ")

# --- Step 2: append a brand-new run after that, holding the red annotation.
$para1 = $d.Paragraphs(1).Range
$noteStart = $para1.End - 1
$noteRange = $d.Range($noteStart, $noteStart)
$noteText = "(This is a change " + [char]0x2013 + " Version for branch alternate)"
$noteRange.InsertAfter($noteText)

$coloredRange = $d.Range($noteStart, $noteStart + $noteText.Length)
$coloredRange.Font.Color = 192
